# Planung.xlsx - edit script
# Applies: window view tweaks, sheet1 selection/scroll changes, re-wording of
# two "Auswertung..." cells into the shorter "Aspektermittlung, Projekt-Setup"
# text, filling the Status column (D) with "OK" for weeks that previously had
# no status, row-height shrink for the now-shorter rows 6/7, and an actuals
# correction on F14 (14.5 -> 6.5) which ripples into the Sum/Average rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Workbook window view ---------------------------------------------------
$wb.Windows.Item(1).WindowState = -4143   # xlNormal, ensure top/left usable
$excel.Windows.Item(1).Top = 45
$excel.Windows.Item(1).Height = 6120

# --- Sheet view: scroll position + selection --------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("F15").Select()

# --- Re-word the "Auswertung bestehende Codebasis..." entries --------------
$ws.Range("C6").Value = "Aspektermittlung, Projekt-Setup"
$ws.Range("C7").Value = "Aspektermittlung, Projekt-Setup"

# Rows 6+7 used to be tall (wrapped, long text) - shrink back down now that
# the text is short.
$ws.Rows.Item(6).RowHeight = 30
$ws.Rows.Item(7).RowHeight = 30

# --- Fill the Status column with "OK" -> copy formatting from a cell that
# already carries the plain "OK" style (C-column / D-column cells all share
# cellXf #1) so we don't mint new style entries.
$fmtSource = $ws.Range("D2")
$fmtSource.Copy()

$okRows = @(7,8,9,10,11,12,13,14,15,16,17,18)
foreach ($r in $okRows) {
    $cell = $ws.Range("D" + $r)
    $cell.Value = "OK"
    $cell.PasteSpecial(-4122)   # xlPasteFormats
}

# D19/D20/D21 stay empty but pick up the same plain style (was the
# left/top-aligned "wrap" style used for the long-text rows).
$emptyRows = @(19,20,21)
foreach ($r in $emptyRows) {
    $cell = $ws.Range("D" + $r)
    $cell.PasteSpecial(-4122)   # xlPasteFormats
}

$excel.CutCopyMode = 0

# --- Actuals correction ------------------------------------------------------
$ws.Range("F14").Value = 6.5

$wb.Application.CalculateFull()
